# Update the "Förändrad" (changed) date column (C) for rows 2-14
# from 2023-09-03 (serial 45172) to 2023-09-06 (serial 45175).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 14; $row++) {
    $ws.Range("C$row").Value = 45175
}
